# "se crean mas signos vitales" - add new vital-sign rows to the resource
# list on Hoja1 (column D), matching the style already used by the
# existing vital-sign entries (D6:D11), and move the active selection
# from E11 to F11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vital signs, one per row, mirroring the formatting of the existing
# "Freciencia Cardiaca" / "PA Sistolica" / "PA Diastolica" entries (D6:D11)
$nuevosSignos = @(
    @{ Row = 12; Text = "temperatura" },
    @{ Row = 13; Text = "Frecuencia Respiratoria" },
    @{ Row = 14; Text = "Saturación O2" },
    @{ Row = 15; Text = "Escala de Dolor (EVA)" }
)

foreach ($signo in $nuevosSignos) {
    $target = $ws.Range("D$($signo.Row)")
    # Copy formatting from an existing vital-sign label cell first so the
    # new cell picks up the same style as D6:D11.
    $ws.Range("D6").Copy($target)
    $target.Value = $signo.Text
}

# Move the active cell/selection to F11 (was E11)
$ws.Range("F11").Select() | Out-Null
